# Swap the presentation's applied theme ("Integral" / Red Violet color
# scheme) for the stock "Office Theme" (Office color scheme) that was
# previously only sitting, unused, as the Notes Master's theme part.
#
# PowerPoint's Theme.ThemeColorScheme exposes the 12 color-scheme slots
# of the slide master's theme (ppt/theme/theme1.xml) in a fixed order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# Re-pointing each slot's RGB to the Office Theme's values effectively
# applies the Office Theme to the deck, matching the OOXML swap.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.Theme.ThemeColorScheme

$scheme.Item(1).RGB  = 0x000000   # dk1
$scheme.Item(2).RGB  = 0xFFFFFF   # lt1
$scheme.Item(3).RGB  = 0x6A5444   # dk2
$scheme.Item(4).RGB  = 0xE6E6E7   # lt2
$scheme.Item(5).RGB  = 0xD59B5B   # accent1
$scheme.Item(6).RGB  = 0x317DED   # accent2
$scheme.Item(7).RGB  = 0xA5A5A5   # accent3
$scheme.Item(8).RGB  = 0x00C0FF   # accent4
$scheme.Item(9).RGB  = 0xC47244   # accent5
$scheme.Item(10).RGB = 0x47AD70   # accent6
$scheme.Item(11).RGB = 0xC16305   # hlink
$scheme.Item(12).RGB = 0x724F95   # folHlink
